$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H53").Value = 464.53845
$ws.Range("I53").Value = 322.23077
$ws.Range("J53").Value = 606.8461
$ws.Range("K53").Value = 322.23077
$ws.Range("L53").Value = 606.8461
$ws.Range("M53").Value = 314.76923
$ws.Range("N53").Value = -1880.8461
$ws.Range("H98").Value = 2534.054
$ws.Range("I98").Value = 890.53845
$ws.Range("J98").Value = 6418.727
$ws.Range("K98").Value = 890.53845
$ws.Range("L98").Value = 6418.727
$ws.Range("M98").Value = 607.46155
$ws.Range("N98").Value = -9414.726999999999
$ws.Range("H118").Value = 747.06665
$ws.Range("J118").Value = 953.6
$ws.Range("L118").Value = 2860.8
$ws.Range("N118").Value = -6174.8
$ws.Range("H122").Value = 2534.054
$ws.Range("I122").Value = 890.53845
$ws.Range("J122").Value = 6418.727
$ws.Range("K122").Value = 2671.61535
$ws.Range("L122").Value = 19256.181
$ws.Range("M122").Value = -221.61535
$ws.Range("N122").Value = -24156.181
$ws.Range("H132").Value = 199616.25
$ws.Range("I132").Value = 3230.3865
$ws.Range("K132").Value = 9691.1595
$ws.Range("M132").Value = -7161.1595
$ws.Range("H135").Value = 312.75
$ws.Range("I135").Value = 263.05554
$ws.Range("J135").Value = 461.83334
$ws.Range("K135").Value = 2367.49986
$ws.Range("L135").Value = 4156.50006
$ws.Range("M135").Value = 167.5001400000001
$ws.Range("N135").Value = -9226.50006
$ws.Range("H137").Value = 2853.6758
$ws.Range("I137").Value = 1474.4783
$ws.Range("J137").Value = 5119.5
$ws.Range("K137").Value = 4423.4349
$ws.Range("L137").Value = 15358.5
$ws.Range("M137").Value = -1873.4349
$ws.Range("N137").Value = -20458.5
$ws.Range("H138").Value = 1996.66
$ws.Range("J138").Value = 2785.6558
$ws.Range("L138").Value = 8356.9674
$ws.Range("N138").Value = -18636.9674

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 1086.4286
$ws.Range("I2").Value = 1012.6667
$ws.Range("J2").Value = 1184.7778
$ws.Range("K2").Value = 1012.6667
$ws.Range("L2").Value = 1184.7778
$ws.Range("M2").Value = -899.6667
$ws.Range("N2").Value = -1410.7778
$ws.Range("H32").Value = 5585.2036
$ws.Range("I32").Value = 4395.326
$ws.Range("J32").Value = 9795.538
$ws.Range("K32").Value = 4395.326
$ws.Range("L32").Value = 9795.538
$ws.Range("M32").Value = -4108.326
$ws.Range("N32").Value = -10369.538
$ws.Range("H35").Value = 20916.25
$ws.Range("J35").Value = 33332.5
$ws.Range("L35").Value = 33332.5
$ws.Range("N35").Value = -34144.5
$ws.Range("H61").Value = 1315.28
$ws.Range("I61").Value = 804.8889
$ws.Range("J61").Value = 2627.7144
$ws.Range("K61").Value = 804.8889
$ws.Range("L61").Value = 2627.7144
$ws.Range("M61").Value = -592.8889
$ws.Range("N61").Value = -3051.7144
$ws.Range("H74").Value = 3431.7646
$ws.Range("I74").Value = 3399.2856
$ws.Range("J74").Value = 3583.3333
$ws.Range("K74").Value = 3399.2856
$ws.Range("L74").Value = 3583.3333
$ws.Range("M74").Value = -2525.2856
$ws.Range("N74").Value = -5331.3333
$ws.Range("H77").Value = 3431.7646
$ws.Range("I77").Value = 3399.2856
$ws.Range("J77").Value = 3583.3333
$ws.Range("K77").Value = 16996.428
$ws.Range("L77").Value = 17916.6665
$ws.Range("M77").Value = -12628.428
$ws.Range("N77").Value = -26652.6665
$ws.Range("H116").Value = 1086.4286
$ws.Range("I116").Value = 1012.6667
$ws.Range("J116").Value = 1184.7778
$ws.Range("K116").Value = 1012.6667
$ws.Range("L116").Value = 1184.7778
$ws.Range("M116").Value = 1281.3333
$ws.Range("N116").Value = -5772.7778
$ws.Range("H136").Value = 1315.28
$ws.Range("I136").Value = 804.8889
$ws.Range("J136").Value = 2627.7144
$ws.Range("K136").Value = 2414.6667
$ws.Range("L136").Value = 7883.1432
$ws.Range("M136").Value = 135.3332999999998
$ws.Range("N136").Value = -12983.1432
$ws.Range("H139").Value = 43034.617
$ws.Range("J139").Value = 43034.617
$ws.Range("L139").Value = 43034.617
$ws.Range("N139").Value = -53314.617

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 1086.4286
$ws.Range("I3").Value = 1012.6667
$ws.Range("J3").Value = 1184.7778
$ws.Range("K3").Value = 1012.6667
$ws.Range("L3").Value = 1184.7778
$ws.Range("M3").Value = -898.6667
$ws.Range("N3").Value = -1412.7778
$ws.Range("H134").Value = 1803.7241
$ws.Range("I134").Value = 1052.762
$ws.Range("J134").Value = 3775
$ws.Range("K134").Value = 3158.286
$ws.Range("L134").Value = 11325
$ws.Range("M134").Value = -623.2860000000001
$ws.Range("N134").Value = -16395
$ws.Range("H138").Value = 41324.445
$ws.Range("J138").Value = 41324.445
$ws.Range("L138").Value = 41324.445
$ws.Range("N138").Value = -51604.445

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 20003228
$ws.Range("I31").Value = 1766.0667
$ws.Range("J31").Value = 50005424
$ws.Range("K31").Value = 1766.0667
$ws.Range("L31").Value = 50005424
$ws.Range("M31").Value = -1471.0667
$ws.Range("N31").Value = -50006014
$ws.Range("H34").Value = 20003228
$ws.Range("I34").Value = 1766.0667
$ws.Range("J34").Value = 50005424
$ws.Range("K34").Value = 1766.0667
$ws.Range("L34").Value = 50005424
$ws.Range("M34").Value = -1564.0667
$ws.Range("N34").Value = -50005828
$ws.Range("H132").Value = 3849.6365
$ws.Range("I132").Value = 1792.5714
$ws.Range("J132").Value = 7449.5
$ws.Range("K132").Value = 5377.7142
$ws.Range("L132").Value = 22348.5
$ws.Range("M132").Value = -2847.7142
$ws.Range("N132").Value = -27408.5
$ws.Range("H134").Value = 7468.4
$ws.Range("I134").Value = 13449.5
$ws.Range("J134").Value = 3481
$ws.Range("K134").Value = 40348.5
$ws.Range("L134").Value = 10443
$ws.Range("M134").Value = -37813.5
$ws.Range("N134").Value = -15513
$ws.Range("H138").Value = 45308.89
$ws.Range("J138").Value = 45308.89
$ws.Range("L138").Value = 45308.89
$ws.Range("N138").Value = -55588.89
$ws.Range("H140").Value = 93212.73
$ws.Range("J140").Value = 93212.73
$ws.Range("L140").Value = 93212.73
$ws.Range("N140").Value = -103572.73
$ws.Range("H141").Value = 34000
$ws.Range("J141").Value = 34000
$ws.Range("L141").Value = 34000
$ws.Range("N141").Value = -44360

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H97").Value = 480.2857
$ws.Range("I97").Value = 227
$ws.Range("J97").Value = 2000
$ws.Range("K97").Value = 681
$ws.Range("L97").Value = 6000
$ws.Range("M97").Value = -185
$ws.Range("N97").Value = -6992
$ws.Range("H113").Value = 612.1053000000001
$ws.Range("I113").Value = 583.53845
$ws.Range("J113").Value = 674
$ws.Range("K113").Value = 1750.61535
$ws.Range("L113").Value = 2022
$ws.Range("M113").Value = 419.38465
$ws.Range("N113").Value = -6362
$ws.Range("H137").Value = 2772.4
$ws.Range("I137").Value = 910
$ws.Range("J137").Value = 3339.2173
$ws.Range("K137").Value = 2730
$ws.Range("L137").Value = 10017.6519
$ws.Range("M137").Value = 2370
$ws.Range("N137").Value = -20217.6519

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H140").Value = 40457.07
$ws.Range("J140").Value = 40457.07
$ws.Range("L140").Value = 40457.07
$ws.Range("N140").Value = -50817.07

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H46").Value = 2066.7778
$ws.Range("I46").Value = 3633.3333
$ws.Range("J46").Value = 1753.4667
$ws.Range("K46").Value = 3633.3333
$ws.Range("L46").Value = 1753.4667
$ws.Range("M46").Value = -3445.3333
$ws.Range("N46").Value = -2129.4667
$ws.Range("H93").Value = 4445743
$ws.Range("I93").Value = 6173593
$ws.Range("J93").Value = 2700
$ws.Range("K93").Value = 6173593
$ws.Range("L93").Value = 2700
$ws.Range("M93").Value = -6172345
$ws.Range("N93").Value = -5196
$ws.Range("H132").Value = 4703.4326
$ws.Range("I132").Value = 1884.375
$ws.Range("J132").Value = 9907.846
$ws.Range("K132").Value = 5653.125
$ws.Range("L132").Value = 29723.538
$ws.Range("M132").Value = -3123.125
$ws.Range("N132").Value = -34783.538
$ws.Range("H139").Value = 41710
$ws.Range("J139").Value = 41710
$ws.Range("L139").Value = 41710
$ws.Range("N139").Value = -51990
$ws.Range("H140").Value = 91428.42999999999
$ws.Range("J140").Value = 91428.42999999999
$ws.Range("L140").Value = 91428.42999999999
$ws.Range("N140").Value = -101788.43
$ws.Range("H141").Value = 40745.625
$ws.Range("J141").Value = 40745.625
$ws.Range("L141").Value = 40745.625
$ws.Range("N141").Value = -51105.625

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H138").Value = 43355.332
$ws.Range("J138").Value = 43355.332
$ws.Range("L138").Value = 43355.332
$ws.Range("N138").Value = -53635.332
$ws.Range("H139").Value = 37653.168
$ws.Range("J139").Value = 37549.83
$ws.Range("L139").Value = 37549.83
$ws.Range("N139").Value = -47829.83
$ws.Range("H140").Value = 52005.4
$ws.Range("J140").Value = 52005.4
$ws.Range("L140").Value = 52005.4
$ws.Range("N140").Value = -62365.4
$ws.Range("H141").Value = 43250
$ws.Range("J141").Value = 43250
$ws.Range("L141").Value = 43250
$ws.Range("N141").Value = -53610
